$d = $word.ActiveDocument

$p7 = $d.Paragraphs.Item(7)
$full7 = $p7.Range
$del7 = $d.Range($full7.Start, $full7.End - 1)
$del7.Delete()
$ins7 = $d.Range($full7.Start, $full7.Start)
$xml7 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Add the ESLint npm package using npm or yarn:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins7.InsertXML($xml7)

$p11 = $d.Paragraphs.Item(11)
$full11 = $p11.Range
$del11 = $d.Range($full11.Start, $full11.End - 1)
$del11.Delete()
$ins11 = $d.Range($full11.Start, $full11.Start)
$xml11 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Setup the initial configuration using the init command and following the steps:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins11.InsertXML($xml11)

$p13 = $d.Paragraphs.Item(13)
$full13 = $p13.Range
$del13 = $d.Range($full13.Start, $full13.End - 1)
$del13.Delete()
$ins13 = $d.Range($full13.Start, $full13.Start)
$xml13 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/></w:r><w:r w:rsidRPr="00404BEE"><w:t>yarn eslint --init</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins13.InsertXML($xml13)

$p32 = $d.Paragraphs.Item(32)
$full32 = $p32.Range
$del32 = $d.Range($full32.Start, $full32.End - 1)
$del32.Delete()
$ins32 = $d.Range($full32.Start, $full32.Start)
$xml32 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>To make linting directories quicker, add the eslint script to the package.json, for example:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins32.InsertXML($xml32)

$p34 = $d.Paragraphs.Item(34)
$full34 = $p34.Range
$del34 = $d.Range($full34.Start, $full34.End - 1)
$del34.Delete()
$ins34 = $d.Range($full34.Start, $full34.Start)
$xml34 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>scripts": {</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins34.InsertXML($xml34)

$p36 = $d.Paragraphs.Item(36)
$full36 = $p36.Range
$del36 = $d.Range($full36.Start, $full36.End - 1)
$del36.Delete()
$ins36 = $d.Range($full36.Start, $full36.Start)
$xml36 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">    "eslint": "eslint --fix \"app/javascript/packs/controllers</w:t></w:r><w:r w:rsidR="008C5005"><w:t>/**/*</w:t></w:r><w:r w:rsidR="009B7E3B"><w:t>\""</w:t></w:r><w:r w:rsidR="008C5005"><w:t>,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins36.InsertXML($xml36)

$p37 = $d.Paragraphs.Item(37)
$full37 = $p37.Range
$del37 = $d.Range($full37.Start, $full37.End - 1)
$del37.Delete()
$ins37 = $d.Range($full37.Start, $full37.Start)
$xml37 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">    "eslint-spec": "eslint --fix \"spec/javascript/**/*\""</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins37.InsertXML($xml37)

$p50 = $d.Paragraphs.Item(50)
$full50 = $p50.Range
$del50 = $d.Range($full50.Start, $full50.End - 1)
$del50.Delete()
$ins50 = $d.Range($full50.Start, $full50.Start)
$xml50 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Then add the parser details to the .eslintrc configuration file:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins50.InsertXML($xml50)

$p55 = $d.Paragraphs.Item(55)
$full55 = $p55.Range
$del55 = $d.Range($full55.Start, $full55.End - 1)
$del55.Delete()
$ins55 = $d.Range($full55.Start, $full55.Start)
$xml55 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>ESLint configuration files are hierarchical and can therefore be added to directories to give them specific configuration. For example, to turn off some rules for a specific directory, add a new .eslintrc file in the directory with the following:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins55.InsertXML($xml55)

$p57 = $d.Paragraphs.Item(57)
$full57 = $p57.Range
$del57 = $d.Range($full57.Start, $full57.End - 1)
$del57.Delete()
$ins57 = $d.Range($full57.Start, $full57.Start)
$xml57 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>{</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins57.InsertXML($xml57)

$p59 = $d.Paragraphs.Item(59)
$full59 = $p59.Range
$del59 = $d.Range($full59.Start, $full59.End - 1)
$del59.Delete()
$ins59 = $d.Range($full59.Start, $full59.Start)
$xml59 = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">    "no-undef": "off"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins59.InsertXML($xml59)

$p61 = $d.Paragraphs.Item(61)
$full61 = $p61.Range
$ins61 = $d.Range($full61.End - 1, $full61.End - 1)
$xmlNew = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t>Ignore</w:t></w:r></w:p><w:p><w:r><w:t>Files and directories can be ignored using a .eslintignore file, similar to .gitignore</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins61.InsertXML($xmlNew)
